$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "changed" date for every data row (2..150).
# All of them are being bumped from serial date 45172 (2023-09-03) to
# serial date 45175 (2023-09-06). Set the underlying numeric value directly
# so the existing date number-format (style index 1) is preserved.
$newDate = [DateTime]::FromOADate(45175)

for ($row = 2; $row -le 150; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
